$d = $word.ActiveDocument

$targetCount = 29

$p = $d.Paragraphs.Item(1)
$p.Range.Text = "⚡️🚀המאמר היומי של מייק -07.12.24: ⚡️🚀" + [char]11 + "Retrieval-Augmented Generation with Knowledge Graphs for Customer Service Question Answering"

$p = $d.Paragraphs.Item(2)
$p.Range.Text = "1. תמצית המאמר"
$p.Style = "Heading4"

$p = $d.Paragraphs.Item(3)
$p.Range.Text = "המאמר מציע שיטה המציידת RAG עם מערכת מבוססת גרפי ידע (KG) המותאמת לשירות לקוחות. המערכת, שפותחה על ידי צוות המחקר של LinkedIn, מעשירה LLMs בידע מבני שמקורו בפניות שירות היסטוריות. על ידי שילוב יחסים שונים בין פניות השירות (טיקטים) בגרף, השיטה משפרת באופן משמעותי את דיוק האחזור(retrieval), איכות התשובות והיעילות, עם שיפורים ניכרים במדדים כמו MRR, BLEU ומקטין זמני הטיפול בפניות."

$p = $d.Paragraphs.Item(4)
$p.Range.Text = "2. תרומות מרכזיות"
$p.Style = "Heading4"

$p = $d.Paragraphs.Item(5)
$p.Range.Text = "א. שילוב KG במערכות RAG"
$p.Style = "Heading5"

$p = $d.Paragraphs.Item(6)
$p.Range.Text = "שימור מידע מבני:" + [char]11 + "כל טיקט מיוצג כעץ (יחסים פנימיים בתוכו) ומקושרת לפניות אחרות דרך יחסים סמנטיים או מפורשים. עיצוב זה משמר את ההיגיון הלוגי של הטיקט, כולל תיאור הבעיה והפתרון. כל טיקט מהווה צומת בגרף."

$p = $d.Paragraphs.Item(7)
$p.Range.Text = "שיפור באחזור ויצירת תשובות:" + [char]11 + "המערכת מנווטת בגרף כדי לזהות תתי-גרפים רלוונטיים, המוזנים ל-LLMs לצורך יצירת תשובות איכותיות."

$p = $d.Paragraphs.Item(8)
$p.Range.Text = "ב. בניית גרף הידע:"
$p.Style = "Heading5"

$p = $d.Paragraphs.Item(9)
$p.Range.Text = "עץ פנימי לטיקט:" + [char]11 + "צמתים מייצגים חלקים כמו סיכומים או שורשי בעיה, וקשתות מציינות יחסים היררכיים."

$p = $d.Paragraphs.Item(10)
$p.Range.Text = "קשרים בין פניות:"

$p = $d.Paragraphs.Item(11)
$p.Range.Text = "קשרים מפורשים: יחסים כמו (e.g., `"clone of`" or `"caused by`")."

$p = $d.Paragraphs.Item(12)
$p.Range.Text = "קשרים סמויים: מחושבים על בסיס דמיון קוסיין בין אמבדינגס."

$p = $d.Paragraphs.Item(13)
$p.Range.Text = "ג. שלבים בתהליך אחזור ותשובות"
$p.Style = "Heading5"

$p = $d.Paragraphs.Item(14)
$p.Range.Text = "המערכת פועלת ב 3 שלבים:"

$p = $d.Paragraphs.Item(15)
$p.Range.Text = "זיהוי ישויות(entity) וכוונות:" + [char]11 + "המערכת הופכת שאילתות משתמש לישויות וכוונות(intents) באמצעות LLMs וניתוח ותבניות YAML."

$p = $d.Paragraphs.Item(16)
$p.Range.Text = "אחזור תת-גרפים:" + [char]11 + "מתבצע חישוב דמיון בין אמבדינגס לשאילתה לצמתים בגרף לזיהוי תת-הגרפים הרלוונטיים ביותר."

$p = $d.Paragraphs.Item(17)
$p.Range.Text = "יצירת תשובות:" + [char]11 + "המערכת יוצרת תשובות בהתבסס על תת-הגרפים רלוונטיים לשאילתת המשתמש."

$p = $d.Paragraphs.Item(18)
$p.Range.Text = "4. קצת פרטים על השיטה"

$p = $d.Paragraphs.Item(19)
$p.Range.Text = "השיטה המוצעת כוללת 3 שלבים עיקריים:"

$p = $d.Paragraphs.Item(20)
$p.Range.Text = "a. זיהוי ישויות בשאילתה וזיהוי כוונה(intent):"

$p = $d.Paragraphs.Item(21)
$p.Range.Text = "המערכת מעבדת שאילתות משתמש על ידי חילוץ ישויות מוגדרות וכוונות באמצעות ניתוח תבניות YAML ו-LLMs. ישויות מוגדרות מייצגות אופיינים מהותיים (למשל, `"תקציר בעיה`" או `"תיאור בעיה`"), בעוד כוונות(intents) מכילות את מטרת השאילתה (למשל, `"פתרון תיקון`"). לדוגמה, בהינתן השאילתה `"כיצד לשחזר את בעיית ההתחברות כאשר משתמש לא יכול להתחבר ל-LinkedIn?`", המערכת מזהה את הישויות כ`"בעיית התחברות`" ו`"משתמש לא יכול להתחבר`" ואת הכוונה כ`"פתרון תיקון.`""

$p = $d.Paragraphs.Item(22)
$p.Range.Text = "b.  אחזור מבוסס אמבדינגס (ייצוג):"

$p = $d.Paragraphs.Item(23)
$p.Range.Text = "זיהוי פניות רלוונטיות: מחשבים עד כמה הישויות שחולצו משאילתת המשתמש (למשל, `"בעיית התחברות`") תואמות את הצמתים ב-KG. עבור כל יישות בשאילתה, השיטה משתמשת בדמיון קוסיין למדידת קרבה בין ייצוג הישות לייצוגים של צמתים בגרף. הציונים מצטברים על פני כל הצמתים השייכים לטיקט מסוים. ככל שלטיקט יש מספר ישויות קרובות לשאילתה, הציון שלו עולה, מה שהופך אותו לסביר יותר להיבחר כרלוונטי. " + [char]11

$p = $d.Paragraphs.Item(24)
$p.Range.Text = "חילוץ תת-גרף רלוונטי: לאחר זיהוי טיקטים הרלוונטיים ביותר, הם משמשים לבניית שאילותות למסד נתונים (DB) בשפת שאילתות גרפים הנקראת Cypher. שאילתות אלה מאפשרות למערכת לחלץ תת-גרפים מקושרים, כגון תיאורים קשורים או שלבים לשחזור בעיה. תהליך האחזור המובנה הזה מבטיח(סוג של כמו תמיד) שהמערכת אוספת מידע מדויק ורלוונטי מבחינת ההקשר מגרף הידע."

$p = $d.Paragraphs.Item(25)
$p.Range.Text = "c. יצירת תשובה:"

$p = $d.Paragraphs.Item(26)
$p.Range.Text = "מגנרטת תשובות על ידי קישור נתוני הגרף שאוחזרו עם השאילתה המקורית. LLM מנסח מחדש את השאילתה באופן דינמי ומייצר תשובות מובנות. לדוגמה השאילתה `"שגיאת העלאת csv בעדכון אימייל משתמש`" מנוסחת מחדש ל-Cypher לאינטראקציה עם DB, מאחזרת פתרונות צעד-אחר-צעד."

$p = $d.Paragraphs.Item(27)
$p.Range.Text = "6. סיכום"
$p.Style = "Heading4"

$p = $d.Paragraphs.Item(28)
$p.Range.Text = "המאמר מציג דרך פורצת דרך לשילוב גרפי ידע במערכות RAG עבור מענה לשאלות בשירות לקוחות. על ידי לכידת יחסים פנימיים וחיצוניים בין פניות, המערכת משפרת משמעותית את דיוק האחזור ואיכות יצירת התשובות, ומציבה כיוון מעניין ביישומים פרקטיים של LLMs."

$p = $d.Paragraphs.Item(29)
$p.Range.Text = "https://arxiv.org/abs/2404.17723:"

if ($d.Paragraphs.Count -gt $targetCount) {
    $startPara = $d.Paragraphs.Item($targetCount + 1)
    $delRange = $d.Range($startPara.Range.Start, $d.Content.End)
    $delRange.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count